$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$newPlaceholder = "XXXXXXXXXXXXXXXXXXXXXXXX"

# Fill in new data dictionary rows 8 and 9, in the same order the values were
# originally typed so shared-string indices line up (column by column).
$ws.Range("B8").Value = "dl"
$ws.Range("B9").Value = "slash"
$ws.Range("C8").Value = "icône téléchargement"
$ws.Range("C9").Value = "icône code slash"
$ws.Range("D8").Value = "elle contient mon cv"
$ws.Range("D9").Value = $newPlaceholder

# Update the old placeholder text (previously 25 X's) in rows 4-6 in column D
# to the newly entered 24 X's placeholder value.
$ws.Range("D4").Value = $newPlaceholder
$ws.Range("D5").Value = $newPlaceholder
$ws.Range("D6").Value = $newPlaceholder

# Match the border/style used by row 8 (thin border all around) for row 9
# by copying the existing cell format instead of creating a brand-new style.
$ws.Range("B8:D8").Copy()
$ws.Range("B9:D9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Update the active selection to E6, matching the saved view state
$ws.Range("E6").Select()
